$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, shifting existing rows 11..98 down to 12..99
$ws.Range("A11").EntireRow.Insert()

# Populate the newly inserted row 11 with the new weekly price record
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Vega Modelo de Temuco"
$ws.Range("C11").Value = "La Araucanía"
$ws.Range("D11").Value = 45230
$ws.Range("E11").Value = 9
$ws.Range("F11").Value = 100112026
$ws.Range("G11").Value = "Haba"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 65
$ws.Range("K11").Value = 12000
$ws.Range("L11").Value = 12000
$ws.Range("M11").Value = 12000
$ws.Range("N11").Value = "$/saco 25 kilos"
$ws.Range("O11").Value = "Región del Maule"
$ws.Range("P11").Value = 480
$ws.Range("Q11").Value = 25
$ws.Range("R11").Value = "Hortaliza"
